# Update "Tasks" workbook: add scoring annotations to the PA1-PA4 headers,
# add a new PA5 column (F) with its own header and per-row notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PA")

# --- Update header row (B1:E1) with score annotations ---
$ws.Range("B1").Value = "PA1 (20/20)"
$ws.Range("C1").Value = "PA2 (100/100)"
$ws.Range("D1").Value = "PA3 (90/100)"
$ws.Range("E1").Value = "PA4 (113+13/120+20)"

# --- New column F: header + per-student notes ---
$ws.Range("F1").Value = "PA5 (?/20+50)"
$ws.Range("F2").Value = "· Tổng hợp câu hỏi khảo sát"
$ws.Range("F3").Value = "· Chọn câu hỏi khảo sát"
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = "· Chụp hình và chú thích`n· Chọn câu hỏi khảo sát"
$ws.Range("F6").Value = "· Chọn câu hỏi khảo sát"

# --- Formatting: copy styles from the neighboring column E so the new
# column F matches the look (borders/fill/alignment/wrap) of the table ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

$ws.Range("E2:E6").Copy() | Out-Null
$ws.Range("F2:F6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column width for the new column (best-effort match of the ~24 char
# width used elsewhere in the sheet) ---
$ws.Columns.Item(6).ColumnWidth = 23.4

# Leave the selection on F5, matching the last-edited cell
$ws.Range("F5").Select()
